$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: C1 50 -> 200; D1 200 -> empty (clear value, keep style)
$ws.Range("C1").Value = 200
$ws.Range("D1").ClearContents()

# Row 2: remove D2 value
$ws.Range("D2").ClearContents()

# Row 3: remove D3 value
$ws.Range("D3").ClearContents()

# Row 4: C4 0.2 -> 0.63; remove D4 value
$ws.Range("C4").Value = 0.63
$ws.Range("D4").ClearContents()

# Row 5: remove D5 value
$ws.Range("D5").ClearContents()

# Row 6: C6 0.5 -> 0.765; remove D6 value
$ws.Range("C6").Value = 0.76500000000000001
$ws.Range("D6").ClearContents()

# Row 7: remove D7 value (formula)
$ws.Range("D7").ClearContents()

# Row 8: remove D8 value
$ws.Range("D8").ClearContents()

# Update selection to C7
$ws.Range("C7").Select()
